$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27

# Column A holds a literal text string "05/22/2025" (not a real date value),
# matching the other rows in this sheet that store dates as text.
# Temporarily force text format so Excel doesn't auto-convert it to a date,
# then restore the default "Normal" style so no stray style survives.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "05/22/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 422.9830000000002
$ws.Cells.Item($row, 3).Value = 0.1182080603712205
$ws.Cells.Item($row, 4).Value = 50
